$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column U (shifts "nom"/"url_produit" from U/V to V/W).
$ws.Columns("U").Insert()

# Populate the new column U with a copy of column T's values (history snapshot column),
# matching the pattern used for every prior timestamp column.
$ws.Range("T1:T205").Copy()
$ws.Range("U1").PasteSpecial(-4163)

# The new snapshot column's header gets today's timestamp instead of a copy of "reference".
$ws.Range("U1").Value = "2026-01-28 14:21:16"
